{"js": "// Replace each \"before\" division fact with its \"after\" counterpart.\n// All values are unique in the document, so a straightforward\n// search-and-replace on each distinct string is safe and order-independent.\nconst replacements = [\n  [\"32\u00f75=6, 2\", \"57\u00f79=6, 3\"],\n  [\"81\u00f72=40, 1\", \"51\u00f72=25, 1\"],\n  [\"51\u00f78=6, 3\", \"80\u00f76=13, 2\"],\n  [\"74\u00f74=18, 2\", \"36\u00f77=5, 1\"],\n  [\"69\u00f79=7, 6\", \"34\u00f75=6, 4\"],\n  [\"95\u00f75=19, 0\", \"21\u00f73=7, 0\"],\n  [\"26\u00f74=6, 2\", \"57\u00f78=7, 1\"],\n  [\"61\u00f78=7, 5\", \"86\u00f75=17, 1\"],\n  [\"88\u00f77=12, 4\", \"85\u00f79=9, 4\"],\n  [\"55\u00f73=18, 1\", \"51\u00f79=5, 6\"],\n  [\"94\u00f73=31, 1\", \"58\u00f73=19, 1\"],\n  [\"12\u00f78=1, 4\", \"23\u00f79=2, 5\"],\n  [\"26\u00f72=13, 0\", \"99\u00f78=12, 3\"],\n  [\"83\u00f78=10, 3\", \"78\u00f74=19, 2\"],\n  [\"67\u00f75=13, 2\", \"76\u00f77=10, 6\"],\n  [\"89\u00f76=14, 5\", \"88\u00f76=14, 4\"],\n  [\"20\u00f76=3, 2\", \"75\u00f74=18, 3\"],\n  [\"96\u00f78=12, 0\", \"52\u00f74=13, 0\"],\n  [\"29\u00f77=4, 1\", \"22\u00f73=7, 1\"],\n  [\"16\u00f75=3, 1\", \"90\u00f78=11, 2\"],\n  [\"57\u00f76=9, 3\", \"21\u00f76=3, 3\"],\n  [\"37\u00f72=18, 1\", \"43\u00f75=8, 3\"],\n  [\"86\u00f76=14, 2\", \"75\u00f79=8, 3\"],\n  [\"15\u00f74=3, 3\", \"16\u00f73=5, 1\"],\n  [\"13\u00f78=1, 5\", \"74\u00f72=37, 0\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"before\" division fact with its \"after\" counterpart.\n# Every original value is unique in the document, so a plain Find/Replace\n# (ReplaceAll, MatchCase on, no wildcards) per pair is safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"32\u00f75=6, 2\", \"57\u00f79=6, 3\"),\n    @(\"81\u00f72=40, 1\", \"51\u00f72=25, 1\"),\n    @(\"51\u00f78=6, 3\", \"80\u00f76=13, 2\"),\n    @(\"74\u00f74=18, 2\", \"36\u00f77=5, 1\"),\n    @(\"69\u00f79=7, 6\", \"34\u00f75=6, 4\"),\n    @(\"95\u00f75=19, 0\", \"21\u00f73=7, 0\"),\n    @(\"26\u00f74=6, 2\", \"57\u00f78=7, 1\"),\n    @(\"61\u00f78=7, 5\", \"86\u00f75=17, 1\"),\n    @(\"88\u00f77=12, 4\", \"85\u00f79=9, 4\"),\n    @(\"55\u00f73=18, 1\", \"51\u00f79=5, 6\"),\n    @(\"94\u00f73=31, 1\", \"58\u00f73=19, 1\"),\n    @(\"12\u00f78=1, 4\", \"23\u00f79=2, 5\"),\n    @(\"26\u00f72=13, 0\", \"99\u00f78=12, 3\"),\n    @(\"83\u00f78=10, 3\", \"78\u00f74=19, 2\"),\n    @(\"67\u00f75=13, 2\", \"76\u00f77=10, 6\"),\n    @(\"89\u00f76=14, 5\", \"88\u00f76=14, 4\"),\n    @(\"20\u00f76=3, 2\", \"75\u00f74=18, 3\"),\n    @(\"96\u00f78=12, 0\", \"52\u00f74=13, 0\"),\n    @(\"29\u00f77=4, 1\", \"22\u00f73=7, 1\"),\n    @(\"16\u00f75=3, 1\", \"90\u00f78=11, 2\"),\n    @(\"57\u00f76=9, 3\", \"21\u00f76=3, 3\"),\n    @(\"37\u00f72=18, 1\", \"43\u00f75=8, 3\"),\n    @(\"86\u00f76=14, 2\", \"75\u00f79=8, 3\"),\n    @(\"15\u00f74=3, 3\", \"16\u00f73=5, 1\"),\n    @(\"13\u00f78=1, 5\", \"74\u00f72=37, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1  # wdFindContinue\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute(\n        $find,        # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replace,     # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n"}
